$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.024.63'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.20%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.909.76'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.66%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.67%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.56'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.08%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4826'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3810'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07358'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9333'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.00%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.921.61'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.52%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.642'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.48%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '92.00'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.95%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.006'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008876'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.83%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '28.040.83'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.74'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.159'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.162.94'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.61%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.89'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '156.98'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.916'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.50'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.118'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +5.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '117.16'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.54%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.987'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.75%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.76%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.284'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.257'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7760'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.83%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.596'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -4.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02058'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.108'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5532'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.43%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05301'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.62%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.001'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.017'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.56%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1529'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.35%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.516'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.73'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '109.18'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +6.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4826'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.005'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.649'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.84%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '68.10'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06077'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.06%  '
